# Updates cryptos list figures (prices / 1h volume %) and reorders the
# Uniswap/Polkadot rows, per the "Updated cryptos list" GitHub Actions commit.
#
# Price/volume cells hold text (not numbers) even when the text looks like a
# plain number (e.g. "1.00"), so for any such D-column cell we briefly force
# the cell to Text format before writing the value (otherwise Excel would
# auto-coerce "1.00" into the number 1 and drop the formatting), then put the
# format back to General so no stray style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.962.25"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "3.388.90"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.40"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.27"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D8").Value = "3.388.57"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.64"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "3.969.59"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.07"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "3.390.24"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "61.130.10"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.99"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.81"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.46"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.10"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.555"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "3.529.36"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.12"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.78"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +11.86%  "
$ws.Range("E29").Value = "  +9.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.52"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.16"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.68"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  -4.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.88"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.64"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0760"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  -4.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.66"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.41"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.29"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "2.467.05"
$ws.Range("E48").Value = "  +3.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.18"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.81"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.43"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +5.20%  "
